$d = $word.ActiveDocument

# --- Add a new paragraph style "AbstractTitle" ("Abstract Title") ---
# It sits right before the existing "Abstract" style in the stylesheet,
# based on Normal, with Abstract as the "next style".
$s = $d.Styles.Add("AbstractTitle", 1)
$s.NameLocal = "Abstract Title"
$s.BaseStyle = "Normal"
$s.NextParagraphStyle = "Abstract"
$s.QuickStyle = $true

$s.ParagraphFormat.KeepWithNext = $true   # w:keepNext
$s.ParagraphFormat.KeepTogether = $true   # w:keepLines
$s.ParagraphFormat.Alignment = 1          # wdAlignParagraphCenter -> w:jc center
$s.ParagraphFormat.SpaceBefore = 15       # 300 twips
$s.ParagraphFormat.SpaceAfter = 0         # 0 twips

$s.Font.Size = 10                         # w:sz 20 (half-points)
$s.Font.SizeBi = 10                       # w:szCs 20
$s.Font.Bold = $true                      # w:b
$s.Font.Color = 9067060                   # RGB(0x34,0x5A,0x8A) -> w:color 345A8A

# --- Tweak the existing "Abstract" style: reduce space-before to 100 ---
$abstract = $d.Styles("Abstract")
$abstract.ParagraphFormat.SpaceBefore = 5  # 100 twips (spacing after stays 300)
